# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
# Rebuilds the "Source" section of the Summary sheet (rows 49-65):
#  - removes the old hyperlink on the source URL cell
#  - inserts new "Additional Information" / "Informal Sector Information" blocks
#  - updates the OEDE description text
#  - re-spaces everything with blank separator rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlink that lived on A51 (the trabajo.gov.ar source link).
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Range("A51").Hyperlinks.Delete()
}

# Clear out the old block (rows 49-56) before laying the new one down.
$ws.Range("A49:E65").Clear()

# Row 49: "Source:" heading
$ws.Range("A49").Value = "Source:"
$ws.Range("A49").Style = "source"

# Row 50: blank separator
$ws.Range("A50").Style = "source"

# Row 51: citation text
$ws.Range("A51").Value = 'Ministerio de Trabajo, Empleo y Seguridad Social. "Caracterización y evolución de la cantidad de empresas. Serie anual."'
$ws.Range("A51").Style = "source"

# Row 52: blank separator
$ws.Range("A52").Style = "source"

# Row 53: source URL (plain text now, no hyperlink)
$ws.Range("A53").Value = "http://www.trabajo.gov.ar/left/estadisticas/oede/estadisticas_nacionales.asp"
$ws.Range("A53").Style = "source"

# Row 54: blank separator
$ws.Range("A54").Style = "source"

# Row 55: "Additional Information:" heading
$ws.Range("A55").Value = "Additional Information:"
$ws.Range("A55").Style = "source"

# Row 56: blank separator
$ws.Range("A56").Style = "source"

# Row 57: ECLAC link
$ws.Range("A57").Value = "http://www.eclac.org/publicaciones/xml/0/11180/lcbuel178.pdf"
$ws.Range("A57").Style = "source"

# Row 58: blank separator
$ws.Range("A58").Style = "source"

# Row 59: "Informal Sector Information:" heading
$ws.Range("A59").Value = "Informal Sector Information:"
$ws.Range("A59").Style = "source"

# Row 60: blank separator
$ws.Range("A60").Style = "source"

# Row 61: World Bank microdata link
$ws.Range("A61").Value = "http://microdata.worldbank.org/index.php/catalog/411"
$ws.Range("A61").Style = "source"

# Rows 62-63 intentionally left blank (no cells)

# Row 64: "OEDE" title
$ws.Range("A64").Value = "OEDE"
$ws.Range("A64").Style = "title"

# Row 65: updated OEDE / SEPyME resolution description
$ws.Range("A65").Value = "Resolution 24/2001 from Secretariat of Small and Medium Enterprise (SEPyME), Ministry of Economy. Modified in 20 August 2010 by Resolution 21/2010 in relation to the determination of the annual sales value."
$ws.Range("A65").Style = "source"
